$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = 51
$ws.Range("B26").Value = "try sq 2"
$ws.Range("C26").Value = "riya-morankar"
$ws.Range("D26").Value = "N/A"
$ws.Range("E26").Value = "edit1 to main"

$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "2025-06-20"
$ws.Range("F26").ClearFormats()
